# Apply the "added scraping code for extra bowling attributes" edit:
#  1. Add a new worksheet "ODI Bowling Extra" after "ODI Batting Extra"
#     with the MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL table.
#  2. On the existing "ODI Batting Extra" sheet, clear out the leftover
#     blank B/C/D/E placeholder cells on rows that have no batting data
#     (match codes 4247, 4276, 4287, 4294, 4337, 4432, 4433), leaving
#     just the MATCH_CODE (A) and MAN_OF_MATCH (F) cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Clean up "ODI Batting Extra" - drop the empty B/C/D/E cells
# ---------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$battingExtra.Range("B2:E2").ClearContents()
$battingExtra.Range("B6:E6").ClearContents()
$battingExtra.Range("B8:E8").ClearContents()
$battingExtra.Range("B10:E10").ClearContents()
$battingExtra.Range("B16:E16").ClearContents()
$battingExtra.Range("B20:E20").ClearContents()
$battingExtra.Range("B21:E21").ClearContents()

# ---------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" worksheet
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Reuse the header formatting (bold + border + centered) from the
# "ODI Batting Extra" header row.
$battingExtra.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Data rows are stored as plain text (same convention as the rest of
# the workbook), so force the Text number format before writing values
# to stop Excel from re-interpreting "0"/"1"/"10.00%" as numbers.
$bowlingExtra.Range("A2:C21").NumberFormat = "@"

$data = @(
  @("4241", "0", ""),
  @("4244", "0", "10.00%"),
  @("4247", "", ""),
  @("4273", "1", ""),
  @("4274", "1", ""),
  @("4275", "0", "10.00%"),
  @("4276", "", ""),
  @("4277", "1", ""),
  @("4292", "0", ""),
  @("4294", "", ""),
  @("4297", "0", "20.00%"),
  @("4300", "0", "30.00%"),
  @("4324", "0", ""),
  @("4334", "0", ""),
  @("4337", "", ""),
  @("4340", "0", "20.00%"),
  @("4349", "0", ""),
  @("4375", "1", "10.00%"),
  @("4432", "", ""),
  @("4433", "", "")
)

$r = 2
foreach ($row in $data) {
  if ($row[0] -ne "") {
    $bowlingExtra.Cells.Item($r, 1).Value = $row[0]
  }
  if ($row[1] -ne "") {
    $bowlingExtra.Cells.Item($r, 2).Value = $row[1]
  }
  if ($row[2] -ne "") {
    $bowlingExtra.Cells.Item($r, 3).Value = $row[2]
  }
  $r = $r + 1
}

$battingExtra.Select()
$battingExtra.Range("A1").Select()
